$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New one-hot "Engine" classification columns E:G for rows 3-12
$data = @(
    @(0, 0, 1),
    @(0, 0, 1),
    @(0, 1, 1),
    @(0, 1, 0),
    @(0, 0, 1),
    @(1, 0, 0),
    @(0, 1, 0),
    @(0, 1, 0),
    @(1, 0, 0),
    @(0, 1, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 3 + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 5).Value = $row[0]
    $ws.Cells.Item($r, 6).Value = $row[1]
    $ws.Cells.Item($r, 7).Value = $row[2]
}

# Totals row 13: SUM formulas for B:D, literal totals for E:G
$ws.Range("B13").Formula = "=SUM(B3:B12)"
$ws.Range("C13").Formula = "=SUM(C3:C12)"
$ws.Range("D13").Formula = "=SUM(D3:D12)"

$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 5
$ws.Range("G13").Value = 4

# View state: zoom level and active selection
$excel.ActiveWindow.Zoom = 244
$ws.Range("F14").Select() | Out-Null
